$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shift existing D:K -> F:M)
$ws.Columns("D:E").Insert()

# --- Re-apply formatting to the full D:M block (rows 7-102) ---
# Date header rows (7, 38, 80): date format + bold font, matches surrounding cells
# (NOTE: multi-area comma ranges only format the first area in this host, so issue
#  one statement per contiguous area instead of a single comma-joined Range.)
$r = $ws.Range("D7:M7")
$r.NumberFormat = "[$-409]d\-mmm\-yy;@"
$r.Font.Bold = $true
$r.Font.Name = "Verdana"
$r.Font.Size = 12
$r = $ws.Range("D38:M38")
$r.NumberFormat = "[$-409]d\-mmm\-yy;@"
$r.Font.Bold = $true
$r.Font.Name = "Verdana"
$r.Font.Size = 12
$r = $ws.Range("D80:M80")
$r.NumberFormat = "[$-409]d\-mmm\-yy;@"
$r.Font.Bold = $true
$r.Font.Name = "Verdana"
$r.Font.Size = 12

# Data rows: #,##0 number format, right aligned, non-bold Verdana 12 (matches existing data cells)
$r = $ws.Range("D8:M35")
$r.NumberFormat = "#,##0"
$r.HorizontalAlignment = -4152
$r.Font.Bold = $false
$r.Font.Name = "Verdana"
$r.Font.Size = 12
$r = $ws.Range("D41:M77")
$r.NumberFormat = "#,##0"
$r.HorizontalAlignment = -4152
$r.Font.Bold = $false
$r.Font.Name = "Verdana"
$r.Font.Size = 12
$r = $ws.Range("D81:M102")
$r.NumberFormat = "#,##0"
$r.HorizontalAlignment = -4152
$r.Font.Bold = $false
$r.Font.Name = "Verdana"
$r.Font.Size = 12

# --- Write cell values/content for D:M across all data rows ---
$rowsData = @{}
$rowsData[7] = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
$rowsData[8] = @(272100, 176100, 192500, 179300, 231400, 187400, 229400, 167500, 207900, 158300)
$rowsData[9] = @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
$rowsData[10] = @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
$rowsData[11] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowsData[12] = @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
$rowsData[13] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[14] = @("NA", "NA", 2500, 0, 0, "NA", 0, 0, -6300, "NA")
$rowsData[15] = @("NA", 49500, 46200, 45700, 52300, 50400, 41500, 42900, 48800, 43700)
$rowsData[16] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowsData[17] = @(214700, 158000, 189200, 182400, 176400, 164000, 175800, 172900, 158300, 150200)
$rowsData[18] = @(57400, 18100, 3300, -3100, 55000, 23400, 53600, -5400, 49600, 8100)
$rowsData[19] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowsData[20] = @(-8300, -200, -11300, 7000, 700, 5600, 3000, 2200, -8700, -7800)
$rowsData[21] = @(96600, 67400, 38200, 49600, 108000, 79400, 98100, 39700, 89700, 44100)
$rowsData[22] = @("NA", 6500, 6900, 5700, 6400, 6700, 7000, 6500, 7100, "NA")
$rowsData[23] = @(49100, 11400, -14900, -1700, 49300, 22300, 49600, -9700, 33800, 300)
$rowsData[24] = @(14200, 2700, 6300, -1800, 11100, 9900, 18000, -4800, 10700, -1200)
$rowsData[25] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[26] = @(34900, 8700, -21200, 100, 38200, 12500, 31600, -5000, 23100, 1500)
$rowsData[27] = @(34900, 8700, -21200, 100, 38200, 12500, 31600, -5000, 23100, 1500)
$rowsData[28] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[29] = @("NA", "NA", 0, 0, -63900, "NA", "NA", "NA", "NA", "NA")
$rowsData[30] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[31] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[32] = @(8300, 200, 11300, -7000, -700, -5600, -3000, -2200, 8700, 7800)
$rowsData[33] = @(34900, 8700, -21200, 100, -25700, 12500, 31600, -5000, 23100, 1500)
$rowsData[34] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[35] = @(34900, 8700, -21200, 100, -25700, 12500, 31600, -5000, 23100, 1500)
$rowsData[38] = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
$rowsData[39] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowsData[40] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowsData[41] = @(248300, 322800, 315700, 210900, 202600, 282000, 378600, 201000, 161600, 216300)
$rowsData[42] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[43] = @(162800, 120900, 220200, 118500, 160500, 111200, 135300, 114100, 163800, 90500)
$rowsData[44] = @(3400, 3500, 1800, 1800, 2500, 2700, 2100, 1800, 1400, 1900)
$rowsData[45] = @(800, 700, 1500, 3700, 3500, 3200, 4200, 5200, 6000, 6900)
$rowsData[46] = @(415400, 447800, 539300, 334900, 369100, 399100, 520200, 322100, 332900, 315500)
$rowsData[47] = @(14300, 13900, 6900, 8700, 16200, 13000, 20100, 19700, 5600, 5300)
$rowsData[48] = @(340000, 340900, 356200, 337900, 661600, 340100, 337400, 339900, 339900, 341100)
$rowsData[49] = @(964500, 1001000, 1043000, 980900, 2008800, 1050900, 935900, 931700, 1018200, 1053800)
$rowsData[50] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[51] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[52] = @(78500, 87300, 88900, 109300, 102600, 177500, 186700, 192900, 194100, 199300)
$rowsData[53] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[54] = @(1812700, 1891000, 2016100, 1771600, 1832200, 1980600, 2000400, 1806400, 1890700, 1915000)
$rowsData[55] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowsData[56] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowsData[57] = @(235500, 241300, 349500, 272400, 238600, 264200, 248200, 232300, 219500, 224800)
$rowsData[58] = @(7200, 3000, 11800, 7800, 11900, 8100, 7500, 3600, 8100, 3500)
$rowsData[59] = @(176600, 296400, 240500, 102800, 157200, 262900, 279700, 141200, 204600, 266400)
$rowsData[60] = @(419400, 540700, 601900, 383000, 407800, 535200, 535400, 377100, 432200, 494800)
$rowsData[61] = @(655500, 642300, 634800, 596100, 619200, 623500, 649100, 679800, 692400, 657500)
$rowsData[62] = @(147300, 144000, 222400, 184000, 181100, 172900, 189700, 146700, 150200, 163300)
$rowsData[63] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[64] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[65] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[66] = @(1222100, 1326900, 1459100, 1163100, 1208100, 1331600, 1374200, 1203600, 1274800, 1315500)
$rowsData[67] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowsData[68] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[69] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[70] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[71] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[72] = @(547300, 512100, 503200, 535200, 565500, 590400, 577100, 562200, 582100, 558300)
$rowsData[73] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[74] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[75] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[76] = @(590500, 564100, 557100, 608500, 624100, 649100, 626200, 602800, 615800, 599500)
$rowsData[77] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[80] = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
$rowsData[81] = @(34900, 8700, -21200, 100, -25700, 12500, 31600, -5000, 23100, 1500)
$rowsData[82] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowsData[83] = @(47500, 49500, 46200, 45700, 52300, 50400, 41500, 42900, 48800, 43700)
$rowsData[84] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[85] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[86] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[87] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[88] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[89] = @(-55300, 149700, 131200, 27600, -57900, 23400, 231800, 52400, -55900, 69300)
$rowsData[90] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowsData[91] = @(-3100, -6400, -4800, -1300, -5500, -5700, -2600, -9900, -11500, -211200)
$rowsData[92] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[93] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[94] = @(-24200, -141700, -21500, 500, -20800, -116400, -31700, 5300, -7700, -164500)
$rowsData[95] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$rowsData[96] = @(0, 0, -14400, 0, 0, 0, -15000, -15600, 0, 0)
$rowsData[97] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[98] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[99] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowsData[100] = @(0, -4900, -14600, -14400, -100, -100, -15100, -15700, -100, -100)
$rowsData[101] = @(5100, 4000, 9600, -5400, -600, -3300, -5600, -2500, 9000, 9800)
$rowsData[102] = @(-74500, 7200, 104700, 8400, -79500, -96600, 179500, 39400, -54700, -85500)

foreach ($r in $rowsData.Keys) {
  $vals = $rowsData[$r]
  for ($i = 0; $i -lt $vals.Count; $i++) {
    $v = $vals[$i]
    if ($null -ne $v) {
      $ws.Cells.Item([int]$r, 4 + $i).Value = $v
    }
  }
}

